$wb = $excel.ActiveWorkbook

# --- Add the new "Worlds-like test" sheet as the last tab, after Sheet7 ---
$sheet7 = $wb.Worksheets.Item("Sheet7")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet7)
$ws.Name = "Worlds-like test"

# --- Header row ---
$ws.Range("A1").Value = "COORD"
$ws.Range("B1").Value = "Predicted X"
$ws.Range("C1").Value = "Predicted Y"
$ws.Range("E1").Value = "Real X (TAG REL)"
$ws.Range("F1").Value = "Real Y"
$ws.Range("K1").Value = "Delta X"
$ws.Range("L1").Value = "Delta Y"
$ws.Range("M1").Value = "Delta H (deg)"

# --- Point labels first (A2:A6), so the shared-string table picks up
#     B1 / B3 / B4 / B5 in that order (B2 already exists in the table) ---
$ws.Range("A2").Value = "B1"
$ws.Range("A3").Value = "B2"
$ws.Range("A4").Value = "B3"
$ws.Range("A5").Value = "B4"
$ws.Range("A6").Value = "B5"

# --- then the "AVG X ERR" / "PROCESSED" labels ---
$ws.Range("K22").Value = "AVG X ERR"
$ws.Range("K10").Value = "PROCESSED"

# --- and only at the end the "SPIKE" callout label ---
$ws.Range("A8").Value = "SPIKE (B5) - important point"

# --- Row 2 : point B1 ---
$ws.Range("B2").Value = 48
$ws.Range("C2").Value = -44
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = 47.5
$ws.Range("F2").Value = 47.625
$ws.Range("K2").Formula = "= E2 - ABS(B2)"
$ws.Range("L2").Formula = "= F2 - ABS(C2)"
$ws.Range("M2").Formula = "= G2 - ABS(D2)"

# --- Row 3 : point B2 ---
$ws.Range("B3").Value = 23.7
$ws.Range("C3").Value = -44.5
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = 23.625
$ws.Range("F3").Value = 47.625
$ws.Range("H3").Formula = "=5/8"
$ws.Range("K3").Formula = "= E3 - ABS(B3)"
$ws.Range("L3").Formula = "= F3 - ABS(C3)"
$ws.Range("M3").Formula = "= G3 - ABS(D3)"

# --- Row 4 : point B3 ---
$ws.Range("B4").Value = 47.4
$ws.Range("C4").Value = -19.9
$ws.Range("D4").Value = $null
$ws.Range("E4").Value = 47.5
$ws.Range("F4").Value = 23.875
$ws.Range("K4").Formula = "= E4 - ABS(B4)"
$ws.Range("L4").Formula = "= F4 - ABS(C4)"
$ws.Range("M4").Formula = "= G4 - ABS(D4)"

# --- Row 5 : point B4 ---
$ws.Range("B5").Value = 23.6
$ws.Range("C5").Value = -20.3
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = 23.75
$ws.Range("F5").Value = 23.875
$ws.Range("K5").Formula = "= E5 - ABS(B5)"
$ws.Range("L5").Formula = "= F5 - ABS(C5)"
$ws.Range("M5").Formula = "= G5 - ABS(D5)"

# --- Row 6 : point B5 (the outlier / spike) ---
$ws.Range("B6").Value = -1.8
$ws.Range("C6").Value = -21.3
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 23.75
$ws.Range("K6").Formula = "= E6 - ABS(B6)"
$ws.Range("L6").Formula = "= F6 - ABS(C6)"
$ws.Range("M6").Formula = "= G6 - ABS(D6)"

# --- Row 8 : SPIKE (B5) - important point ---
$ws.Range("B8").Value = 14.4
$ws.Range("C8").Value = -20.6
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = 15.625
$ws.Range("F8").Value = 23.375
$ws.Range("K8").Formula = "= E8 - ABS(B8)"
$ws.Range("L8").Formula = "= F8 - ABS(C8)"
$ws.Range("M8").Formula = "= G8 - ABS(D8)"

# --- Row 12 : comparison of the spike row vs the averages ---
$ws.Range("K12").Formula = "=K8-K23"
$ws.Range("L12").Formula = "=L8-L23"

# --- Averages ---
$ws.Range("L22").Value = "AVG Y ERR"
$ws.Range("K23").Formula = "=AVERAGE(K2:K5, K8)"
$ws.Range("L23").Formula = "=AVERAGE(L2:L5, L8)"

# --- Deviations from the average, one per calibration point ---
$ws.Range("K25").Formula = "=K2-`$K`$23"
$ws.Range("L25").Formula = "=L2-`$L`$23"
$ws.Range("K26").Formula = "=K3-`$K`$23"
$ws.Range("L26").Formula = "=L3-`$L`$23"
$ws.Range("K27").Formula = "=K4-`$K`$23"
$ws.Range("L27").Formula = "=L4-`$L`$23"
$ws.Range("K28").Formula = "=K5-`$K`$23"
$ws.Range("L28").Formula = "=L5-`$L`$23"
